$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy formatting (bold/border/alignment)
# from the neighboring header cell G1 so H1 gets the exact same style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column in row 2
$ws.Range("H2").Value = 1
